# Loan RBI, Variable Instalments
#
# Insert a new (blank) column before the existing "Late" column on the
# "Repayment schedule" sheet so a Variable-Instalment column can be added,
# pushing Late / Outstanding(heading) / Outstanding one column to the right.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Column N ("Late") is the 14th column - insert a new blank column there,
# shifting N->O, O->P, P->Q (and widening the sheet from P13 to Q13).
$ws.Columns.Item(14).Insert()

# Match the new column's width to its neighbour (column M, width 11) so the
# stored width comes out exactly "11" instead of a rounded character-width
# conversion.
$ws.Columns.Item(14).ColumnWidth = $ws.Columns.Item(13).ColumnWidth

# Reflect the final on-screen state: the "Repayment schedule" tab active,
# with cell I18 selected.
$ws.Activate()
$ws.Range("I18").Select()
